$d = $word.ActiveDocument

# The "Proposed Improvement" section has a run of empty paragraphs right
# after its heading. The 2nd of those empty paragraphs gains a hanging
# indent plus the new body text (as two runs), and the following two
# empty paragraphs are removed outright.
$target = $d.Paragraphs.Item(41)

# Add the paragraph indent (w:ind w:left="360" w:firstLine="720") -
# LeftIndent/FirstLineIndent are expressed in points; 18pt == 360 twips,
# 36pt == 720 twips.
$target.LeftIndent = 18
$target.FirstLineIndent = 36

# Fill in the paragraph text, matching the existing paragraph-mark run
# formatting (Arial, 12pt / 24 half-points, incl. complex-script fields).
$sentence = "The researchers encountered issues in this topic, such as the absence of a scientific calculator that saves the user's recent calculations. Therefore, the researchers proposed the following improvement to reduce the hassle for users: exploring different approaches, such as using a database engine to store recent calculations"

$r = $target.Range
$r.Text = $sentence + "."
$r.Font.Name = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Size = 12
$r.Font.SizeBi = 12

# Split the trailing "." into its own run (matches the source markup,
# which has the sentence and the final period as two separate <w:r>
# elements) by nudging formatting on just that character and back.
$paraStart = $target.Range.Start
$periodStart = $paraStart + $sentence.Length
$rPeriod = $d.Range($periodStart, $periodStart + 1)
$rPeriod.Bold = 1
$rPeriod.Bold = 0

# Remove the next two (now-empty) paragraphs entirely.
$d.Paragraphs.Item(42).Range.Delete()
$d.Paragraphs.Item(42).Range.Delete()
